$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-02-01 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-02 Sunday", 2) | Out-Null
$d.Content.Find.Execute("61-26=35", $true, $false, $false, $false, $false, $true, 1, $false, "49+45=94", 2) | Out-Null
$d.Content.Find.Execute("90-87=3", $true, $false, $false, $false, $false, $true, 1, $false, "91-67=24", 2) | Out-Null
$d.Content.Find.Execute("43+18=61", $true, $false, $false, $false, $false, $true, 1, $false, "46+28=74", 2) | Out-Null
$d.Content.Find.Execute("6+16=22", $true, $false, $false, $false, $false, $true, 1, $false, "24-19=5", 2) | Out-Null
$d.Content.Find.Execute("37+25=62", $true, $false, $false, $false, $false, $true, 1, $false, "50-15=35", 2) | Out-Null
$d.Content.Find.Execute("8+36=44", $true, $false, $false, $false, $false, $true, 1, $false, "52-33=19", 2) | Out-Null
$d.Content.Find.Execute("52-18=34", $true, $false, $false, $false, $false, $true, 1, $false, "75-46=29", 2) | Out-Null
$d.Content.Find.Execute("90-76=14", $true, $false, $false, $false, $false, $true, 1, $false, "60-9=51", 2) | Out-Null
$d.Content.Find.Execute("41-35=6", $true, $false, $false, $false, $false, $true, 1, $false, "34-29=5", 2) | Out-Null
$d.Content.Find.Execute("42-39=3", $true, $false, $false, $false, $false, $true, 1, $false, "38+53=91", 2) | Out-Null
$d.Content.Find.Execute("25+56=81", $true, $false, $false, $false, $false, $true, 1, $false, "59+23=82", 2) | Out-Null
$d.Content.Find.Execute("64-7=57", $true, $false, $false, $false, $false, $true, 1, $false, "32-23=9", 2) | Out-Null
$d.Content.Find.Execute("95-69=26", $true, $false, $false, $false, $false, $true, 1, $false, "77+6=83", 2) | Out-Null
$d.Content.Find.Execute("61-56=5", $true, $false, $false, $false, $false, $true, 1, $false, "18+56=74", 2) | Out-Null
$d.Content.Find.Execute("59+34=93", $true, $false, $false, $false, $false, $true, 1, $false, "84-27=57", 2) | Out-Null
$d.Content.Find.Execute("36+35=71", $true, $false, $false, $false, $false, $true, 1, $false, "13+8=21", 2) | Out-Null
$d.Content.Find.Execute("22-16=6", $true, $false, $false, $false, $false, $true, 1, $false, "91-13=78", 2) | Out-Null
$d.Content.Find.Execute("58+25=83", $true, $false, $false, $false, $false, $true, 1, $false, "71-66=5", 2) | Out-Null
$d.Content.Find.Execute("29+39=68", $true, $false, $false, $false, $false, $true, 1, $false, "70-31=39", 2) | Out-Null
$d.Content.Find.Execute("81-2=79", $true, $false, $false, $false, $false, $true, 1, $false, "29+66=95", 2) | Out-Null
$d.Content.Find.Execute("43-19=24", $true, $false, $false, $false, $false, $true, 1, $false, "39+4=43", 2) | Out-Null
$d.Content.Find.Execute("30-12=18", $true, $false, $false, $false, $false, $true, 1, $false, "81-55=26", 2) | Out-Null
$d.Content.Find.Execute("67+5=72", $true, $false, $false, $false, $false, $true, 1, $false, "6+55=61", 2) | Out-Null
$d.Content.Find.Execute("63-18=45", $true, $false, $false, $false, $false, $true, 1, $false, "43-18=25", 2) | Out-Null
$d.Content.Find.Execute("70-64=6", $true, $false, $false, $false, $false, $true, 1, $false, "87-79=8", 2) | Out-Null
$d.Content.Find.Execute("44-39=5", $true, $false, $false, $false, $false, $true, 1, $false, "19+8=27", 2) | Out-Null
$d.Content.Find.Execute("26-7=19", $true, $false, $false, $false, $false, $true, 1, $false, "9+84=93", 2) | Out-Null
$d.Content.Find.Execute("35+59=94", $true, $false, $false, $false, $false, $true, 1, $false, "17+26=43", 2) | Out-Null
$d.Content.Find.Execute("66+28=94", $true, $false, $false, $false, $false, $true, 1, $false, "35+16=51", 2) | Out-Null
$d.Content.Find.Execute("52-8=44", $true, $false, $false, $false, $false, $true, 1, $false, "56+39=95", 2) | Out-Null
$d.Content.Find.Execute("10-3=7", $true, $false, $false, $false, $false, $true, 1, $false, "19+74=93", 2) | Out-Null
$d.Content.Find.Execute("95-49=46", $true, $false, $false, $false, $false, $true, 1, $false, "65-27=38", 2) | Out-Null
$d.Content.Find.Execute("21-4=17", $true, $false, $false, $false, $false, $true, 1, $false, "88+7=95", 2) | Out-Null
$d.Content.Find.Execute("42-33=9", $true, $false, $false, $false, $false, $true, 1, $false, "42-29=13", 2) | Out-Null
$d.Content.Find.Execute("91-27=64", $true, $false, $false, $false, $false, $true, 1, $false, "80-75=5", 2) | Out-Null
$d.Content.Find.Execute("39+16=55", $true, $false, $false, $false, $false, $true, 1, $false, "93-76=17", 2) | Out-Null
$d.Content.Find.Execute("42-16=26", $true, $false, $false, $false, $false, $true, 1, $false, "64-38=26", 2) | Out-Null
$d.Content.Find.Execute("51-2=49", $true, $false, $false, $false, $false, $true, 1, $false, "14+19=33", 2) | Out-Null
$d.Content.Find.Execute("89+5=94", $true, $false, $false, $false, $false, $true, 1, $false, "33-17=16", 2) | Out-Null
$d.Content.Find.Execute("7+84=91", $true, $false, $false, $false, $false, $true, 1, $false, "92-35=57", 2) | Out-Null
$d.Content.Find.Execute("24+7=31", $true, $false, $false, $false, $false, $true, 1, $false, "4+87=91", 2) | Out-Null
$d.Content.Find.Execute("28+16=44", $true, $false, $false, $false, $false, $true, 1, $false, "14+18=32", 2) | Out-Null
$d.Content.Find.Execute("84-29=55", $true, $false, $false, $false, $false, $true, 1, $false, "57+14=71", 2) | Out-Null
$d.Content.Find.Execute("9+13=22", $true, $false, $false, $false, $false, $true, 1, $false, "18+48=66", 2) | Out-Null
$d.Content.Find.Execute("87+4=91", $true, $false, $false, $false, $false, $true, 1, $false, "62-33=29", 2) | Out-Null
$d.Content.Find.Execute("55+16=71", $true, $false, $false, $false, $false, $true, 1, $false, "48+23=71", 2) | Out-Null
$d.Content.Find.Execute("90-78=12", $true, $false, $false, $false, $false, $true, 1, $false, "3+18=21", 2) | Out-Null
$d.Content.Find.Execute("63-39=24", $true, $false, $false, $false, $false, $true, 1, $false, "43-29=14", 2) | Out-Null
$d.Content.Find.Execute("12-9=3", $true, $false, $false, $false, $false, $true, 1, $false, "65-48=17", 2) | Out-Null
$d.Content.Find.Execute("97-59=38", $true, $false, $false, $false, $false, $true, 1, $false, "39+38=77", 2) | Out-Null
$d.Content.Find.Execute("93-88=5", $true, $false, $false, $false, $false, $true, 1, $false, "57-39=18", 2) | Out-Null
$d.Content.Find.Execute("42+49=91", $true, $false, $false, $false, $false, $true, 1, $false, "30-6=24", 2) | Out-Null
$d.Content.Find.Execute("44-17=27", $true, $false, $false, $false, $false, $true, 1, $false, "94-8=86", 2) | Out-Null
$d.Content.Find.Execute("21-9=12", $true, $false, $false, $false, $false, $true, 1, $false, "5+36=41", 2) | Out-Null
$d.Content.Find.Execute("83-17=66", $true, $false, $false, $false, $false, $true, 1, $false, "36+19=55", 2) | Out-Null
$d.Content.Find.Execute("36+26=62", $true, $false, $false, $false, $false, $true, 1, $false, "92-23=69", 2) | Out-Null
$d.Content.Find.Execute("87-49=38", $true, $false, $false, $false, $false, $true, 1, $false, "44+37=81", 2) | Out-Null
$d.Content.Find.Execute("97-29=68", $true, $false, $false, $false, $false, $true, 1, $false, "49+42=91", 2) | Out-Null
$d.Content.Find.Execute("9+65=74", $true, $false, $false, $false, $false, $true, 1, $false, "59+33=92", 2) | Out-Null
$d.Content.Find.Execute("39+23=62", $true, $false, $false, $false, $false, $true, 1, $false, "15-7=8", 2) | Out-Null
$d.Content.Find.Execute("70-51=19", $true, $false, $false, $false, $false, $true, 1, $false, "29+65=94", 2) | Out-Null
$d.Content.Find.Execute("16+46=62", $true, $false, $false, $false, $false, $true, 1, $false, "85-9=76", 2) | Out-Null
$d.Content.Find.Execute("37+27=64", $true, $false, $false, $false, $false, $true, 1, $false, "27+9=36", 2) | Out-Null
$d.Content.Find.Execute("87+6=93", $true, $false, $false, $false, $false, $true, 1, $false, "40-16=24", 2) | Out-Null
$d.Content.Find.Execute("9+7=16", $true, $false, $false, $false, $false, $true, 1, $false, "37+36=73", 2) | Out-Null
$d.Content.Find.Execute("16+6=22", $true, $false, $false, $false, $false, $true, 1, $false, "76-37=39", 2) | Out-Null
$d.Content.Find.Execute("63-56=7", $true, $false, $false, $false, $false, $true, 1, $false, "71-7=64", 2) | Out-Null
$d.Content.Find.Execute("24+68=92", $true, $false, $false, $false, $false, $true, 1, $false, "46+16=62", 2) | Out-Null
$d.Content.Find.Execute("43-26=17", $true, $false, $false, $false, $false, $true, 1, $false, "30-26=4", 2) | Out-Null
$d.Content.Find.Execute("83-74=9", $true, $false, $false, $false, $false, $true, 1, $false, "46+46=92", 2) | Out-Null
$d.Content.Find.Execute("66-48=18", $true, $false, $false, $false, $false, $true, 1, $false, "27+69=96", 2) | Out-Null
$d.Content.Find.Execute("85+6=91", $true, $false, $false, $false, $false, $true, 1, $false, "47+19=66", 2) | Out-Null
$d.Content.Find.Execute("37+17=54", $true, $false, $false, $false, $false, $true, 1, $false, "65+6=71", 2) | Out-Null
$d.Content.Find.Execute("17+65=82", $true, $false, $false, $false, $false, $true, 1, $false, "8+25=33", 2) | Out-Null
$d.Content.Find.Execute("45-9=36", $true, $false, $false, $false, $false, $true, 1, $false, "44-36=8", 2) | Out-Null
$d.Content.Find.Execute("80-62=18", $true, $false, $false, $false, $false, $true, 1, $false, "27+8=35", 2) | Out-Null
$d.Content.Find.Execute("4+78=82", $true, $false, $false, $false, $false, $true, 1, $false, "80-61=19", 2) | Out-Null
$d.Content.Find.Execute("18+24=42", $true, $false, $false, $false, $false, $true, 1, $false, "39+43=82", 2) | Out-Null
$d.Content.Find.Execute("7+14=21", $true, $false, $false, $false, $false, $true, 1, $false, "79+3=82", 2) | Out-Null
$d.Content.Find.Execute("18+73=91", $true, $false, $false, $false, $false, $true, 1, $false, "38+3=41", 2) | Out-Null
$d.Content.Find.Execute("75-16=59", $true, $false, $false, $false, $false, $true, 1, $false, "41-17=24", 2) | Out-Null
$d.Content.Find.Execute("74-16=58", $true, $false, $false, $false, $false, $true, 1, $false, "81-42=39", 2) | Out-Null
$d.Content.Find.Execute("4+17=21", $true, $false, $false, $false, $false, $true, 1, $false, "10-1=9", 2) | Out-Null
$d.Content.Find.Execute("83-29=54", $true, $false, $false, $false, $false, $true, 1, $false, "61-44=17", 2) | Out-Null
$d.Content.Find.Execute("39+22=61", $true, $false, $false, $false, $false, $true, 1, $false, "84-55=29", 2) | Out-Null
$d.Content.Find.Execute("14+48=62", $true, $false, $false, $false, $false, $true, 1, $false, "93-35=58", 2) | Out-Null
$d.Content.Find.Execute("91-12=79", $true, $false, $false, $false, $false, $true, 1, $false, "14+39=53", 2) | Out-Null
$d.Content.Find.Execute("29+22=51", $true, $false, $false, $false, $false, $true, 1, $false, "24-5=19", 2) | Out-Null
$d.Content.Find.Execute("65+29=94", $true, $false, $false, $false, $false, $true, 1, $false, "18+16=34", 2) | Out-Null
$d.Content.Find.Execute("62-58=4", $true, $false, $false, $false, $false, $true, 1, $false, "96-39=57", 2) | Out-Null
$d.Content.Find.Execute("18+46=64", $true, $false, $false, $false, $false, $true, 1, $false, "62-25=37", 2) | Out-Null
$d.Content.Find.Execute("98-89=9", $true, $false, $false, $false, $false, $true, 1, $false, "41-7=34", 2) | Out-Null
$d.Content.Find.Execute("33-16=17", $true, $false, $false, $false, $false, $true, 1, $false, "55-27=28", 2) | Out-Null
$d.Content.Find.Execute("21-14=7", $true, $false, $false, $false, $false, $true, 1, $false, "60-25=35", 2) | Out-Null
$d.Content.Find.Execute("16+27=43", $true, $false, $false, $false, $false, $true, 1, $false, "5+57=62", 2) | Out-Null
$d.Content.Find.Execute("85-38=47", $true, $false, $false, $false, $false, $true, 1, $false, "18+76=94", 2) | Out-Null
$d.Content.Find.Execute("27+55=82", $true, $false, $false, $false, $false, $true, 1, $false, "48+8=56", 2) | Out-Null
$d.Content.Find.Execute("92-5=87", $true, $false, $false, $false, $false, $true, 1, $false, "80-25=55", 2) | Out-Null
$d.Content.Find.Execute("97-88=9", $true, $false, $false, $false, $false, $true, 1, $false, "5+76=81", 2) | Out-Null
$d.Content.Find.Execute("35+46=81", $true, $false, $false, $false, $false, $true, 1, $false, "15+9=24", 2) | Out-Null
